$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / data rows 1-4 (shift values, new header row) ---
$ws.Range("A1").Value = "RowName"
$ws.Range("B1").Value = "Column1"
$ws.Range("C1").Value = "Column2"

$ws.Range("A2").Value = "FirstRow"
$ws.Range("B2").Value = "SecondCell"
$ws.Range("C2").Value = "ThirdCell"

$ws.Range("A3").Value = "SecondRow"
$ws.Range("B3").Value = "SecondCell_1"
$ws.Range("C3").Value = "ThirdCell_1"
$ws.Range("E3").Value = "Khatod"

$ws.Range("A4").Value = "ThirdRow"
$ws.Range("B4").Value = "SecondCell_2"
$ws.Range("C4").Value = "ThirdCell_2"

# --- Rows 5-7 keep Abhi_0 / Abhi_1 / Abhi_2 (already correct values) ---
$ws.Range("A5").Value = "Abhi_0"
$ws.Range("B5").Value = "Abhi_1"
$ws.Range("C5").Value = "Abhi_2"

$ws.Range("A6").Value = "Abhi_0"
$ws.Range("B6").Value = "Abhi_1"
$ws.Range("C6").Value = "Abhi_2"

$ws.Range("A7").Value = "Abhi_0"
$ws.Range("B7").Value = "Abhi_1"
$ws.Range("C7").Value = "Abhi_2"

# --- Remove old rows 8-17 entirely (clear contents) ---
$ws.Range("A8:D17").Clear()

# --- Update selection to match new target ---
$ws.Range("A4").Select()
